$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (was row 4's data)
$ws.Range("D2").Value = 44672
$ws.Range("M2").Value = 8
$ws.Range("O2").Value = 180000
$ws.Range("P2").Value = 180000
$ws.Range("S2").Value = 180000

# Row 3 (was row 2's data)
$ws.Range("D3").Value = 44253
$ws.Range("M3").Value = 12

# Row 4 (was row 3's data)
$ws.Range("D4").Value = 44993
$ws.Range("M4").Value = 14
$ws.Range("O4").Value = 200000
$ws.Range("P4").Value = 190000
$ws.Range("S4").Value = 190000
